$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for changed rows.
# D-column values are forced to text (NumberFormat "@" then restore the
# cell style to Normal) because several of them are decimal-looking
# strings ("534.17", "0.110", ...) that COM would otherwise silently
# coerce into numeric cells, losing the original text formatting.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.477.36"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.03%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.137.60"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.39%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "534.17"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.21%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.25"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.37%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.134.96"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.42%  "
$ws.Range("E9").Value = "  +2.29%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.16"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.21%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.110"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.84%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.395"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.08%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.676.94"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.40%  "
$ws.Range("E14").Value = "  +3.36%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.68"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.44%  "
$ws.Range("E16").Value = "  +0.28%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "58.485.64"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.01%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.153.02"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.42%  "
$ws.Range("E19").Value = "  +0.40%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.90"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.11%  "
$ws.Range("E21").Value = "  -0.96%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "342.90"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.60%  "
$ws.Range("E23").Value = "  +0.22%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.514"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.94%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "67.84"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.81%  "
$ws.Range("E26").Value = "  -0.34%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.11%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0932"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.99%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.53"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.75%  "
$ws.Range("E30").Value = "  -2.54%  "
$ws.Range("E32").Value = "  +1.46%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "21.16"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.06%  "
$ws.Range("E34").Value = "  -0.31%  "
$ws.Range("E37").Value = "  +3.28%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "26.21"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.48%  "
$ws.Range("E39").Value = "  -3.38%  "
$ws.Range("E40").Value = "  +11.81%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0673"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.54%  "
$ws.Range("E42").Value = "  +4.70%  "
$ws.Range("E43").Value = "  +3.74%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.174.75"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.20%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "36.59"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.46%  "
$ws.Range("E46").Value = "  +0.10%  "
$ws.Range("E47").Value = "  +3.09%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.305.90"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.40%  "
$ws.Range("E49").Value = "  +4.60%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "20.69"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.79%  "
$ws.Range("E51").Value = "  +1.87%  "

# Rows 35 and 36 swap places in the ranking (Monero <-> NEARProtocol)
$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.80"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.98%  "
$ws.Range("B36").Value = "Monero"
$ws.Range("C36").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "157.95"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.50%  "
